# 丽水-漫展信息.xlsx update
# - "展览" sheet (index 1): bump two "想去人数" counters and fix the
#   R动漫嘉年华 row's end-date / interest-count.
# - "全部类型" sheet (index 4): same two counter bumps, plus insert the
#   R动漫嘉年华 event as a new row 7 (pushing the LZ栗子 row down to row 8).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: 展览
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item(1)

$ws1.Range("F5").Value = 82
$ws1.Range("F6").Value = 695
$ws1.Range("E7").Value = "2024.08.24 09:30-08.24 17:00"
$ws1.Range("F7").Value = 1

# ---------------------------------------------------------------------
# Sheet 4: 全部类型
# ---------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item(4)

$ws4.Range("F5").Value = 82
$ws4.Range("F6").Value = 695

# Insert a new row 7 (shifts the existing row 7 "LZ栗子" event down to row 8)
$ws4.Rows("7:7").Insert()

$ws4.Range("A7").Value = 6
$ws4.Range("A7").Font.Bold = $true
$ws4.Range("A7").HorizontalAlignment = -4108
$ws4.Range("A7").VerticalAlignment = -4160
$ws4.Range("A7").Borders.LineStyle = 1

$ws4.Range("B7").Value = "'2024-08-24"
$ws4.Range("C7").Value = "丽水·R动漫嘉年华"
$ws4.Range("D7").Value = "中东路848号(解放街交汇) 飞达国际大酒店"
$ws4.Range("E7").Value = "2024.08.24 09:30-08.24 17:00"
$ws4.Range("F7").Value = 1
$ws4.Range("G7").Value = 45
$ws4.Range("H7").Value = "https://show.bilibili.com/platform/detail.html?id=89651"
$ws4.Range("I7").Value = "//i0.hdslb.com/bfs/openplatform/202407/7o5ALbAM1721383424201.jpeg"

# A8 keeps its original numbering (7) after the shift-down
$ws4.Range("A8").Value = 7
